$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRpUNL")

# Recalibrate the retirement-fraction values for these technologies from 0.03 to 0.01
$ws.Range("B2").Value = 0.01   # hard coal
$ws.Range("B3").Value = 0.01   # natural gas steam turbine
$ws.Range("B4").Value = 0.01   # natural gas combined cycle
$ws.Range("B5").Value = 0.01   # nuclear
$ws.Range("B7").Value = 0.01   # onshore wind
$ws.Range("B8").Value = 0.01   # solar PV
$ws.Range("B13").Value = 0.01  # natural gas peaker
$ws.Range("B14").Value = 0.01  # lignite
$ws.Range("B15").Value = 0.01  # offshore wind

# Update the active selection on the sheet to match the saved view state
$prevActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("D11").Select()
$prevActive.Activate()
